$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(28, 8).Value = 346.30768
$ws.Cells.Item(28, 9).Value = 345.81818
$ws.Cells.Item(28, 10).Value = 349
$ws.Cells.Item(28, 11).Value = 345.81818
$ws.Cells.Item(28, 12).Value = 349
$ws.Cells.Item(28, 13).Value = 139.18182
$ws.Cells.Item(28, 14).Value = -1319

$ws.Cells.Item(98, 8).Value = 1565.75
$ws.Cells.Item(98, 10).Value = 1450
$ws.Cells.Item(98, 12).Value = 1450
$ws.Cells.Item(98, 14).Value = -4446

$ws.Cells.Item(107, 8).Value = 715.381
$ws.Cells.Item(107, 9).Value = 598.5294
$ws.Cells.Item(107, 10).Value = 1212
$ws.Cells.Item(107, 11).Value = 598.5294
$ws.Cells.Item(107, 12).Value = 1212
$ws.Cells.Item(107, 13).Value = 1321.4706
$ws.Cells.Item(107, 14).Value = -5052

$ws.Cells.Item(111, 8).Value = 90912110
$ws.Cells.Item(111, 9).Value = 166671840
$ws.Cells.Item(111, 11).Value = 500015520
$ws.Cells.Item(111, 13).Value = -500012453

$ws.Cells.Item(113, 8).Value = 1929.5652
$ws.Cells.Item(113, 9).Value = 1573.75
$ws.Cells.Item(113, 10).Value = 2004.4736
$ws.Cells.Item(113, 11).Value = 1573.75
$ws.Cells.Item(113, 12).Value = 2004.4736
$ws.Cells.Item(113, 13).Value = 1680.25
$ws.Cells.Item(113, 14).Value = -8512.473599999999

$ws.Cells.Item(116, 8).Value = 2027534.9
$ws.Cells.Item(116, 9).Value = 6995175.5
$ws.Cells.Item(116, 10).Value = 3681.3333
$ws.Cells.Item(116, 11).Value = 6995175.5
$ws.Cells.Item(116, 12).Value = 3681.3333
$ws.Cells.Item(116, 13).Value = -6991733.5
$ws.Cells.Item(116, 14).Value = -10565.3333

$ws.Cells.Item(118, 8).Value = 1083.625
$ws.Cells.Item(118, 9).Value = 517.25
$ws.Cells.Item(118, 10).Value = 1650
$ws.Cells.Item(118, 11).Value = 1551.75
$ws.Cells.Item(118, 12).Value = 4950
$ws.Cells.Item(118, 13).Value = 105.25
$ws.Cells.Item(118, 14).Value = -8264

$ws.Cells.Item(122, 8).Value = 1565.75
$ws.Cells.Item(122, 10).Value = 1450
$ws.Cells.Item(122, 12).Value = 4350
$ws.Cells.Item(122, 14).Value = -9250

$ws.Cells.Item(127, 8).Value = 931.5357
$ws.Cells.Item(127, 9).Value = 367.6
$ws.Cells.Item(127, 10).Value = 1054.1305
$ws.Cells.Item(127, 11).Value = 1102.8
$ws.Cells.Item(127, 12).Value = 3162.3915
$ws.Cells.Item(127, 13).Value = 3857.2
$ws.Cells.Item(127, 14).Value = -13082.3915

$ws.Cells.Item(129, 8).Value = 1086.07
$ws.Cells.Item(129, 9).Value = 521.75
$ws.Cells.Item(129, 10).Value = 1109.5834
$ws.Cells.Item(129, 11).Value = 1565.25
$ws.Cells.Item(129, 12).Value = 3328.7502
$ws.Cells.Item(129, 13).Value = 3434.75
$ws.Cells.Item(129, 14).Value = -13328.7502

$ws.Cells.Item(137, 8).Value = 43700.793
$ws.Cells.Item(137, 9).Value = 1624.4
$ws.Cells.Item(137, 10).Value = 113828.11
$ws.Cells.Item(137, 11).Value = 4873.200000000001
$ws.Cells.Item(137, 12).Value = 341484.33
$ws.Cells.Item(137, 13).Value = -2323.200000000001
$ws.Cells.Item(137, 14).Value = -346584.33

$ws.Cells.Item(141, 8).Value = 1998.8572
$ws.Cells.Item(141, 9).Value = 1999.5
$ws.Cells.Item(141, 11).Value = 5998.5
$ws.Cells.Item(141, 13).Value = -818.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(22, 8).Value = 50001000
$ws.Cells.Item(22, 9).Value = 50001000
$ws.Cells.Item(22, 11).Value = 50001000
$ws.Cells.Item(22, 13).Value = -50000701

$ws.Cells.Item(25, 8).Value = 1000
$ws.Cells.Item(25, 9).Value = 1000
$ws.Cells.Item(25, 10).Value = 0
$ws.Cells.Item(25, 11).Value = 1000
$ws.Cells.Item(25, 12).Value = 0
$ws.Cells.Item(25, 13).Value = -598

$ws.Cells.Item(32, 8).Value = 15628313
$ws.Cells.Item(32, 9).Value = 16951394
$ws.Cells.Item(32, 10).Value = 15958.8
$ws.Cells.Item(32, 11).Value = 16951394
$ws.Cells.Item(32, 12).Value = 15958.8
$ws.Cells.Item(32, 13).Value = -16951107
$ws.Cells.Item(32, 14).Value = -16532.8

$ws.Cells.Item(110, 8).Value = 883.7895
$ws.Cells.Item(110, 9).Value = 627.86206
$ws.Cells.Item(110, 10).Value = 1708.4445
$ws.Cells.Item(110, 11).Value = 627.86206
$ws.Cells.Item(110, 12).Value = 1708.4445
$ws.Cells.Item(110, 13).Value = 1417.13794
$ws.Cells.Item(110, 14).Value = -5798.4445

$ws.Cells.Item(122, 8).Value = 1626.5
$ws.Cells.Item(122, 9).Value = 1626.5
$ws.Cells.Item(122, 10).Value = 0
$ws.Cells.Item(122, 11).Value = 4879.5
$ws.Cells.Item(122, 12).Value = 0
$ws.Cells.Item(122, 13).Value = -2429.5

$ws.Cells.Item(132, 8).Value = 1487.3846
$ws.Cells.Item(132, 9).Value = 1144.6086
$ws.Cells.Item(132, 10).Value = 4115.3335
$ws.Cells.Item(132, 11).Value = 3433.8258
$ws.Cells.Item(132, 12).Value = 12346.0005
$ws.Cells.Item(132, 13).Value = -903.8258000000001
$ws.Cells.Item(132, 14).Value = -17406.0005

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 2850.5
$ws.Cells.Item(86, 9).Value = 2280.8
$ws.Cells.Item(86, 10).Value = 3800
$ws.Cells.Item(86, 11).Value = 2280.8
$ws.Cells.Item(86, 12).Value = 3800
$ws.Cells.Item(86, 13).Value = -1157.8
$ws.Cells.Item(86, 14).Value = -6046

$ws.Cells.Item(89, 8).Value = 2850.5
$ws.Cells.Item(89, 9).Value = 2280.8
$ws.Cells.Item(89, 10).Value = 3800
$ws.Cells.Item(89, 11).Value = 11404
$ws.Cells.Item(89, 12).Value = 19000
$ws.Cells.Item(89, 13).Value = -5788
$ws.Cells.Item(89, 14).Value = -30232

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(15, 8).Value = 3901
$ws.Cells.Item(15, 9).Value = 3241.6
$ws.Cells.Item(15, 11).Value = 3241.6
$ws.Cells.Item(15, 13).Value = -3071.6

$ws.Cells.Item(88, 8).Value = 34500
$ws.Cells.Item(88, 10).Value = 34500
$ws.Cells.Item(88, 12).Value = 34500
$ws.Cells.Item(88, 14).Value = -35312

$ws.Cells.Item(91, 8).Value = 34500
$ws.Cells.Item(91, 10).Value = 34500
$ws.Cells.Item(91, 12).Value = 34500
$ws.Cells.Item(91, 14).Value = -37308

$ws.Cells.Item(132, 8).Value = 3017
$ws.Cells.Item(132, 9).Value = 2096.6667
$ws.Cells.Item(132, 10).Value = 4078.923
$ws.Cells.Item(132, 11).Value = 6290.000100000001
$ws.Cells.Item(132, 12).Value = 12236.769
$ws.Cells.Item(132, 13).Value = -3760.000100000001
$ws.Cells.Item(132, 14).Value = -17296.769

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(134, 8).Value = 4791.5
$ws.Cells.Item(134, 9).Value = 1971.6666
$ws.Cells.Item(134, 11).Value = 5914.9998
$ws.Cells.Item(134, 13).Value = -844.9997999999996

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(23, 8).Value = 1000
$ws.Cells.Item(23, 9).Value = 1000
$ws.Cells.Item(23, 11).Value = 1000
$ws.Cells.Item(23, 13).Value = -777

$ws.Cells.Item(110, 8).Value = 34701.6
$ws.Cells.Item(110, 10).Value = 34701.6
$ws.Cells.Item(110, 12).Value = 34701.6
$ws.Cells.Item(110, 14).Value = -42881.6

$ws.Cells.Item(132, 8).Value = 3572.5625
$ws.Cells.Item(132, 9).Value = 3213.6365
$ws.Cells.Item(132, 10).Value = 4362.2
$ws.Cells.Item(132, 11).Value = 9640.9095
$ws.Cells.Item(132, 12).Value = 13086.6
$ws.Cells.Item(132, 13).Value = -7110.9095
$ws.Cells.Item(132, 14).Value = -18146.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 700
$ws.Cells.Item(22, 9).Value = 700
$ws.Cells.Item(22, 11).Value = 700
$ws.Cells.Item(22, 13).Value = -405

$ws.Cells.Item(27, 8).Value = 700
$ws.Cells.Item(27, 9).Value = 700
$ws.Cells.Item(27, 11).Value = 700
$ws.Cells.Item(27, 13).Value = -593

$ws.Cells.Item(30, 8).Value = 768
$ws.Cells.Item(30, 9).Value = 24
$ws.Cells.Item(30, 10).Value = 3000
$ws.Cells.Item(30, 11).Value = 24
$ws.Cells.Item(30, 12).Value = 3000
$ws.Cells.Item(30, 13).Value = 84
$ws.Cells.Item(30, 14).Value = -3216

$ws.Cells.Item(112, 8).Value = 43540
$ws.Cells.Item(112, 10).Value = 43540
$ws.Cells.Item(112, 12).Value = 43540
$ws.Cells.Item(112, 14).Value = -46494

$ws.Cells.Item(132, 8).Value = 8389.475
$ws.Cells.Item(132, 9).Value = 9460.531000000001
$ws.Cells.Item(132, 10).Value = 4105.25
$ws.Cells.Item(132, 11).Value = 28381.593
$ws.Cells.Item(132, 12).Value = 12315.75
$ws.Cells.Item(132, 13).Value = -25851.593
$ws.Cells.Item(132, 14).Value = -17375.75

$ws.Cells.Item(136, 8).Value = 2605.7727
$ws.Cells.Item(136, 9).Value = 1973.3572
$ws.Cells.Item(136, 10).Value = 3712.5
$ws.Cells.Item(136, 11).Value = 5920.071599999999
$ws.Cells.Item(136, 12).Value = 11137.5
$ws.Cells.Item(136, 13).Value = -3370.071599999999
$ws.Cells.Item(136, 14).Value = -16237.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(113, 8).Value = 326.5909
$ws.Cells.Item(113, 9).Value = 294.5238
$ws.Cells.Item(113, 11).Value = 883.5714
$ws.Cells.Item(113, 13).Value = 1286.4286

$ws.Cells.Item(132, 8).Value = 2549.6428
$ws.Cells.Item(132, 9).Value = 2249.182
$ws.Cells.Item(132, 10).Value = 3651.3333
$ws.Cells.Item(132, 11).Value = 6747.545999999999
$ws.Cells.Item(132, 12).Value = 10953.9999
$ws.Cells.Item(132, 13).Value = -4217.545999999999
$ws.Cells.Item(132, 14).Value = -16013.9999

$ws.Cells.Item(136, 8).Value = 3917.628
$ws.Cells.Item(136, 9).Value = 3980.9707
$ws.Cells.Item(136, 10).Value = 3678.3333
$ws.Cells.Item(136, 11).Value = 11942.9121
$ws.Cells.Item(136, 12).Value = 11034.9999
$ws.Cells.Item(136, 13).Value = -9392.9121
$ws.Cells.Item(136, 14).Value = -16134.9999
